$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.39
$ws.Range("E2").Value = 1.36
$ws.Range("B3").Value = 1.43
$ws.Range("G3").Value = 0.6
$ws.Range("B5").Value = 1.36
$ws.Range("D6").Value = 1.57
$ws.Range("E6").Value = 1.34
$ws.Range("G6").Value = 1.01
$ws.Range("C7").Value = 2.27
$ws.Range("F7").Value = 1.52
